# iRipple DTR report: "added colors to rows"
#
# 1. Highlight rows 5-8 and 11-14 (columns A-J) with a red background to
#    flag days that are absent / missing time entries, and mark the
#    "NO OF OVERTIME HOURS" (I) column with 1 for those rows.
# 2. A stray boolean flag gets written into the (merged, hidden) B19 cell.
# 3. Clean up the FLOOR(...) formulas that had a redundant extra ",1" argument.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Color the flagged attendance rows red and set the overtime flag to 1
# ---------------------------------------------------------------------
$redColor = 6184671   # RGB(223, 94, 94) -> 0xDF5E5E, stored as BGR OLE color

$rowGroups = @("A5:J8", "A11:J14")
foreach ($addr in $rowGroups) {
    $rng = $ws.Range($addr)
    $rng.WrapText = $true
    $rng.Interior.Color = $redColor
}

$ws.Range("I5").Value = 1
$ws.Range("I6").Value = 1
$ws.Range("I7").Value = 1
$ws.Range("I8").Value = 1
$ws.Range("I11").Value = 1
$ws.Range("I12").Value = 1
$ws.Range("I13").Value = 1
$ws.Range("I14").Value = 1

# ---------------------------------------------------------------------
# 2. B19 is part of the merged A19:G19 cell, so it has to be unmerged
#    briefly to poke the hidden boolean value into it.
# ---------------------------------------------------------------------
$ws.Range("A19:G19").UnMerge()
$ws.Range("B19").Value = $false
$ws.Range("A19:G19").Merge()
$ws.Range("A19:G19").WrapText = $true

# ---------------------------------------------------------------------
# 3. Simplify the FLOOR(...) formulas (FLOOR(x,1,1) -> FLOOR(x,1))
# ---------------------------------------------------------------------
$ws.Range("B22").Formula = '=FLOOR(F17/8,1)&"."&FLOOR(MOD(F17,8),1)&"."&(MOD(F17,8)-FLOOR(MOD(F17,8),1))*60'
$ws.Range("B23").Formula = '=FLOOR(H19,1)&"."&(H19-FLOOR(H19,1))*8&".0"'
$ws.Range("B24").Formula = '=FLOOR(I19,1)&"."&(I19-FLOOR(I19,1))*8&".0"'
$ws.Range("B27").Formula = '=FLOOR(K27/8,1)&"."&FLOOR(MOD(K27,8),1)&"."&(MOD(K27,8)-FLOOR(MOD(K27,8),1))*60'
